# "payment page without UI"
# Adds a new "RoomImage" field to the data dictionary's RoomType table
# (sheet "แผ่น1"): a new row 39 plus the 3 new shared strings it needs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 already carries the "closing row of a group" formatting
# (thick-ish bottom border: A/C/D/E/F/G/H -> style 6, B -> style 7) that the
# new last row of the RoomType group (row 39) should use, so clone its
# formatting onto row 39 before filling in the new field's data.
$ws.Range("A6:H6").Copy()
$ws.Range("A39:H39").PasteSpecial(-4122)

# New field: RoomImage / Img / (no key) / NULL=N / "รูปตัวอย่างห้อง"
$ws.Range("B39").Value = 6
$ws.Range("C39").Value = "RoomImage"
$ws.Range("E39").Value = "Img"
$ws.Range("G39").Value = "N"
$ws.Range("H39").Value = "รูปตัวอย่างห้อง"

# New row uses the sheet's default row height (15.75) explicitly, matching
# how Excel marks freshly added rows as custom-height.
$ws.Rows.Item(39).RowHeight = 15.75

# Matches the author's last selection before saving.
$ws.Range("I22").Select()
